$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 176.15384
$ws.Range("I2").Value = 86.625
$ws.Range("J2").Value = 319.4
$ws.Range("K2").Value = 86.625
$ws.Range("L2").Value = 319.4
$ws.Range("M2").Value = 26.375
$ws.Range("N2").Value = -545.4
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("H53").Value = 35.833332
$ws.Range("I53").Value = 47.6
$ws.Range("K53").Value = 47.6
$ws.Range("M53").Value = 589.4
$ws.Range("H74").Value = 3811.9443
$ws.Range("I74").Value = 3336.5
$ws.Range("J74").Value = 5476
$ws.Range("K74").Value = 3336.5
$ws.Range("L74").Value = 5476
$ws.Range("M74").Value = -2400.5
$ws.Range("N74").Value = -7348
$ws.Range("H77").Value = 3811.9443
$ws.Range("I77").Value = 3336.5
$ws.Range("J77").Value = 5476
$ws.Range("K77").Value = 16682.5
$ws.Range("L77").Value = 27380
$ws.Range("M77").Value = -12002.5
$ws.Range("N77").Value = -36740
$ws.Range("H100").Value = 3287.4
$ws.Range("I100").Value = 2999
$ws.Range("K100").Value = 2999
$ws.Range("M100").Value = -2458
$ws.Range("H129").Value = 2019.1666
$ws.Range("I129").Value = 1823.4
$ws.Range("J129").Value = 2998
$ws.Range("K129").Value = 5470.200000000001
$ws.Range("L129").Value = 8994
$ws.Range("M129").Value = -470.2000000000007
$ws.Range("N129").Value = -18994
$ws.Range("H133").Value = 80510.875
$ws.Range("J133").Value = 80510.875
$ws.Range("L133").Value = 80510.875
$ws.Range("N133").Value = -90630.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 43784.75
$ws.Range("I74").Value = 78048.69500000001
$ws.Range("J74").Value = 3291
$ws.Range("K74").Value = 78048.69500000001
$ws.Range("L74").Value = 3291
$ws.Range("M74").Value = -77174.69500000001
$ws.Range("N74").Value = -5039
$ws.Range("H77").Value = 43784.75
$ws.Range("I77").Value = 78048.69500000001
$ws.Range("J77").Value = 3291
$ws.Range("K77").Value = 390243.475
$ws.Range("L77").Value = 16455
$ws.Range("M77").Value = -385875.475
$ws.Range("N77").Value = -25191
$ws.Range("H102").Value = 90085.30499999999
$ws.Range("I102").Value = 101784.8
$ws.Range("J102").Value = 51087
$ws.Range("K102").Value = 101784.8
$ws.Range("L102").Value = 51087
$ws.Range("M102").Value = -100162.8
$ws.Range("N102").Value = -54331
$ws.Range("H122").Value = 7600.8335
$ws.Range("J122").Value = 3750
$ws.Range("L122").Value = 11250
$ws.Range("N122").Value = -16150
$ws.Range("H132").Value = 2146.926
$ws.Range("I132").Value = 1624.35
$ws.Range("J132").Value = 3640
$ws.Range("K132").Value = 4873.049999999999
$ws.Range("L132").Value = 10920
$ws.Range("M132").Value = -2343.049999999999
$ws.Range("N132").Value = -15980
$ws.Range("H134").Value = 144999
$ws.Range("J134").Value = 144999
$ws.Range("L134").Value = 144999
$ws.Range("N134").Value = -155139
$ws.Range("H135").Value = 89497.60000000001
$ws.Range("J135").Value = 89497.60000000001
$ws.Range("L135").Value = 89497.60000000001
$ws.Range("N135").Value = -99637.60000000001
$ws.Range("H138").Value = 94994.5
$ws.Range("J138").Value = 94994.5
$ws.Range("L138").Value = 94994.5
$ws.Range("N138").Value = -105274.5
$ws.Range("H139").Value = 90714.664
$ws.Range("J139").Value = 90714.664
$ws.Range("L139").Value = 90714.664
$ws.Range("N139").Value = -100994.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").ClearContents()
$ws.Range("N92").Value = 0
$ws.Range("H99").Value = 3127274.8
$ws.Range("I99").Value = 2214
$ws.Range("J99").Value = 10419083
$ws.Range("K99").Value = 2214
$ws.Range("L99").Value = 10419083
$ws.Range("M99").Value = -716
$ws.Range("N99").Value = -10422079
$ws.Range("H103").Value = 14149.5
$ws.Range("J103").Value = 14149.5
$ws.Range("L103").Value = 14149.5
$ws.Range("N103").Value = -16493.5
$ws.Range("H134").Value = 2275.2083
$ws.Range("J134").Value = 3838.111
$ws.Range("L134").Value = 11514.333
$ws.Range("N134").Value = -16584.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3969.25
$ws.Range("I62").Value = 3994
$ws.Range("J62").Value = 3944.5
$ws.Range("K62").Value = 3994
$ws.Range("L62").Value = 3944.5
$ws.Range("M62").Value = -3370
$ws.Range("N62").Value = -5192.5
$ws.Range("H65").Value = 3969.25
$ws.Range("I65").Value = 3994
$ws.Range("J65").Value = 3944.5
$ws.Range("K65").Value = 19970
$ws.Range("L65").Value = 19722.5
$ws.Range("M65").Value = -16850
$ws.Range("N65").Value = -25962.5
$ws.Range("H99").Value = 3909178.5
$ws.Range("J99").Value = 5211249
$ws.Range("L99").Value = 5211249
$ws.Range("N99").Value = -5214245
$ws.Range("H107").Value = 1109.4
$ws.Range("I107").Value = 1109.4
$ws.Range("K107").Value = 1109.4
$ws.Range("M107").Value = 810.5999999999999
$ws.Range("H126").Value = 3909178.5
$ws.Range("J126").Value = 5211249
$ws.Range("L126").Value = 15633747
$ws.Range("N126").Value = -15638687
$ws.Range("H132").Value = 1754.5834
$ws.Range("I132").Value = 1732.2727
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 5196.8181
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -2666.8181
$ws.Range("N132").Value = -11060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 6.6666665
$ws.Range("I12").Value = 13
$ws.Range("J12").Value = 3.5
$ws.Range("K12").Value = 39
$ws.Range("L12").Value = 10.5
$ws.Range("M12").Value = 134
$ws.Range("N12").Value = -356.5
$ws.Range("H69").Value = 9006
$ws.Range("J69").Value = 8000
$ws.Range("L69").Value = 24000
$ws.Range("N69").Value = -25622
$ws.Range("H72").Value = 9006
$ws.Range("J72").Value = 8000
$ws.Range("L72").Value = 72000
$ws.Range("N72").Value = -80112

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3138.4167
$ws.Range("I126").Value = 1918.5
$ws.Range("J126").Value = 4358.3335
$ws.Range("K126").Value = 5755.5
$ws.Range("L126").Value = 13075.0005
$ws.Range("M126").Value = -3285.5
$ws.Range("N126").Value = -18015.0005
$ws.Range("H132").Value = 5160.3125
$ws.Range("I132").Value = 6143.364
$ws.Range("J132").Value = 2997.6
$ws.Range("K132").Value = 18430.092
$ws.Range("L132").Value = 8992.799999999999
$ws.Range("M132").Value = -15900.092
$ws.Range("N132").Value = -14052.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 411226.47
$ws.Range("I2").Value = 2005620
$ws.Range("J2").Value = 48864.316
$ws.Range("K2").Value = 2005620
$ws.Range("L2").Value = 48864.316
$ws.Range("M2").Value = -2005508
$ws.Range("N2").Value = -49088.316
$ws.Range("H132").Value = 16729.4
$ws.Range("I132").Value = 35825
$ws.Range("J132").Value = 3999
$ws.Range("K132").Value = 107475
$ws.Range("L132").Value = 11997
$ws.Range("M132").Value = -104945
$ws.Range("N132").Value = -17057
$ws.Range("H141").Value = 88558.22
$ws.Range("J141").Value = 71860.57000000001
$ws.Range("L141").Value = 71860.57000000001
$ws.Range("N141").Value = -82220.57000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 145380.78
$ws.Range("J46").Value = 145380.78
$ws.Range("L46").Value = 145380.78
$ws.Range("N46").Value = -145842.78
$ws.Range("H96").Value = 1656.8667
$ws.Range("I96").Value = 1235.3
$ws.Range("J96").Value = 2500
$ws.Range("K96").Value = 1235.3
$ws.Range("L96").Value = 2500
$ws.Range("M96").Value = 137.7
$ws.Range("N96").Value = -5246
$ws.Range("H100").Value = 4203067.5
$ws.Range("I100").Value = 7144579
$ws.Range("J100").Value = 908.4286
$ws.Range("K100").Value = 14289158
$ws.Range("L100").Value = 1816.8572
$ws.Range("M100").Value = -14288617
$ws.Range("N100").Value = -2898.8572
$ws.Range("H107").Value = 6983.324
$ws.Range("I107").Value = 4549.2666
$ws.Range("K107").Value = 13647.7998
$ws.Range("M107").Value = -11727.7998
$ws.Range("H132").Value = 3346539
$ws.Range("I132").Value = 1641.6666
$ws.Range("J132").Value = 6213594
$ws.Range("K132").Value = 4924.9998
$ws.Range("L132").Value = 18640782
$ws.Range("M132").Value = -2394.9998
$ws.Range("N132").Value = -18645842
$ws.Range("H134").Value = 145380.78
$ws.Range("J134").Value = 145380.78
$ws.Range("L134").Value = 436142.34
$ws.Range("N134").Value = -441212.34
$ws.Range("H138").Value = 108997
$ws.Range("J138").Value = 108997
$ws.Range("L138").Value = 108997
$ws.Range("N138").Value = -119277
$ws.Range("H141").Value = 78966.336
$ws.Range("J141").Value = 78966.336
$ws.Range("L141").Value = 78966.336
$ws.Range("N141").Value = -89326.336
